$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D: add "duration" values (rows 3-9), center/vertical aligned + wrap text ---
$durations = @{
    3 = "2h"
    4 = "1d"
    5 = "2h"
    6 = "1d"
    7 = "1d"
    8 = "2d"
    9 = "1h"
}

foreach ($row in 3..9) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.Value = $durations[$row]
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
    $cell.WrapText = $true
}

# --- Column E: update "referrer" text values for rows 4-9 (and add new ones for rows 8-9) ---
$ws.Range("E4").Value = "The HR referrer for the time manger app"
$ws.Range("E5").Value = "All planer department"
$ws.Range("E6").Value = "The planer referrer for the time manager app"
$ws.Range("E7").Value = "All IT department "
$ws.Range("E8").Value = "The IT egenieer referrer for the time manager app"
$ws.Range("E9").Value = "All employees "

# --- Update selection to C3 ---
$ws.Range("C3").Select()
